# Time Log.xlsx - add Sept 16/17 2014 entries, fix totals/percentages.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Set the "Interruption" (D) column first for the new rows - the shared
# formula in column E re-evaluates more reliably when the minutes value
# lands before the Start/Stop times that flip the row from blank to full.
$ws.Range("D46").Value = 5
$ws.Range("D47").Value = 20
$ws.Range("D48").Value = 20
$ws.Range("D49").Value = 0

# Row 45's date moves forward a day (previously 9/15/2014, now 9/16/2014).
$ws.Range("A45").Value = 41898

# Row 46: 9/16/2014, 9:36 AM - 10:51 AM, 5 min interruption, Testing.
$ws.Range("A46").Value = 41898
$ws.Range("B46").Value = 0.39999999999999997
$ws.Range("C46").Value = 0.45208333333333334
$ws.Range("F46").Value = "Testing"

# Row 47: 9/16/2014, 5:18 PM - 6:07 PM, 20 min interruption, Testing.
$ws.Range("A47").Value = 41898
$ws.Range("B47").Value = 0.72083333333333333
$ws.Range("C47").Value = 0.75486111111111109
$ws.Range("F47").Value = "Testing"

# Row 48: 9/16/2014, 11:24 PM - 12:53 AM, 20 min interruption, Coding.
$ws.Range("A48").Value = 41898
$ws.Range("B48").Value = 0.97569444444444453
$ws.Range("C48").Value = 1.0368055555555555
$ws.Range("F48").Value = "Coding"

# Row 49: 9/17/2014, 1:52 PM - 2:27 PM, no interruption, Testing.
$ws.Range("A49").Value = 41899
$ws.Range("B49").Value = 0.57777777777777783
$ws.Range("B49").NumberFormat = "h:mm AM/PM"
$ws.Range("C49").Value = 0.6020833333333333
$ws.Range("F49").Value = "Testing"

# Match the author's final selection/scroll position near the new rows.
$excel.ActiveWindow.TopLeftCell = $ws.Range("A34")
$ws.Range("D50").Select()
